$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, which carries
# the bold/centered/bordered header style) onto the three new header cells
# so they pick up the same style index instead of a freshly-synthesized one.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-41) gets the same season record: 82 wins, 80 losses, 0 ties.
$ws.Range("AD2:AD41").Value = 82
$ws.Range("AE2:AE41").Value = 80
$ws.Range("AF2:AF41").Value = 0
